# Fruta / hortaliza, semanal
# Insert two new weekly price records at the top of the data block (rows 480-481),
# pushing the existing rows 480-539 down to 482-541.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 480:539 down by two rows.
$ws.Rows("480:481").Insert()

# --- New row 480: Acelga, "Primera" quality, week of 2023-08-04 ---
$ws.Cells.Item(480, 1).Value = 7
$ws.Cells.Item(480, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(480, 3).Value = "Ñuble"
$ws.Cells.Item(480, 4).Value = 45142
$ws.Cells.Item(480, 5).Value = 16
$ws.Cells.Item(480, 6).Value = 100112009
$ws.Cells.Item(480, 7).Value = "Acelga"
$ws.Cells.Item(480, 8).Value = "Sin especificar"
$ws.Cells.Item(480, 9).Value = "Primera"
$ws.Cells.Item(480, 10).Value = 150
$ws.Cells.Item(480, 11).Value = 600
$ws.Cells.Item(480, 12).Value = 600
$ws.Cells.Item(480, 13).Value = 600
$ws.Cells.Item(480, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(480, 15).Value = "Región de Ñuble"
$ws.Cells.Item(480, 16).Value = 600
$ws.Cells.Item(480, 17).Value = 1
$ws.Cells.Item(480, 18).Value = "Hortaliza"

# --- New row 481: Acelga, "Segunda" quality, week of 2023-08-04 ---
$ws.Cells.Item(481, 1).Value = 7
$ws.Cells.Item(481, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(481, 3).Value = "Ñuble"
$ws.Cells.Item(481, 4).Value = 45142
$ws.Cells.Item(481, 5).Value = 16
$ws.Cells.Item(481, 6).Value = 100112009
$ws.Cells.Item(481, 7).Value = "Acelga"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Segunda"
$ws.Cells.Item(481, 10).Value = 150
$ws.Cells.Item(481, 11).Value = 500
$ws.Cells.Item(481, 12).Value = 500
$ws.Cells.Item(481, 13).Value = 500
$ws.Cells.Item(481, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(481, 15).Value = "Región de Ñuble"
$ws.Cells.Item(481, 16).Value = 500
$ws.Cells.Item(481, 17).Value = 1
$ws.Cells.Item(481, 18).Value = "Hortaliza"

# --- New rows 540 & 541: duplicate the final two weekly records so the
#     previously-last rows (now at 538/539 after the shift) still appear once
#     more at the tail, matching the published dataset's append pattern.
$ws.Cells.Item(540, 1).Value = 7
$ws.Cells.Item(540, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(540, 3).Value = "Ñuble"
$ws.Cells.Item(540, 4).Value = 44616
$ws.Cells.Item(540, 5).Value = 16
$ws.Cells.Item(540, 6).Value = 100112009
$ws.Cells.Item(540, 7).Value = "Acelga"
$ws.Cells.Item(540, 8).Value = "Sin especificar"
$ws.Cells.Item(540, 9).Value = "Primera"
$ws.Cells.Item(540, 10).Value = 100
$ws.Cells.Item(540, 11).Value = 400
$ws.Cells.Item(540, 12).Value = 450
$ws.Cells.Item(540, 13).Value = 425
$ws.Cells.Item(540, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(540, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(540, 16).Value = 425
$ws.Cells.Item(540, 17).Value = 1
$ws.Cells.Item(540, 18).Value = "Hortaliza"

$ws.Cells.Item(541, 1).Value = 7
$ws.Cells.Item(541, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(541, 3).Value = "Ñuble"
$ws.Cells.Item(541, 4).Value = 44594
$ws.Cells.Item(541, 5).Value = 16
$ws.Cells.Item(541, 6).Value = 100112009
$ws.Cells.Item(541, 7).Value = "Acelga"
$ws.Cells.Item(541, 8).Value = "Sin especificar"
$ws.Cells.Item(541, 9).Value = "Primera"
$ws.Cells.Item(541, 10).Value = 100
$ws.Cells.Item(541, 11).Value = 400
$ws.Cells.Item(541, 12).Value = 450
$ws.Cells.Item(541, 13).Value = 425
$ws.Cells.Item(541, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(541, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(541, 16).Value = 425
$ws.Cells.Item(541, 17).Value = 1
$ws.Cells.Item(541, 18).Value = "Hortaliza"

# Apply the same date-number-format used throughout column D to the
# new date cells (style index 2 in this workbook).
$ws.Cells.Item(480, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(481, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(540, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(541, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Output "done"
